$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: fill in the new BOM line (LED1) ---------------------------------
# Copy cell formatting first (before the values are overwritten) so the
# existing style entries are reused instead of new ones being minted:
#  - old E18 (border box + fill flag) -> A18, C18, D18
#  - old A1  (border box, header style) -> B18, E18, F18
$ws.Range("E18").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("E18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("E18").Copy()
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("F18").PasteSpecial(-4122)

# Now write the new BOM row values
$ws.Range("A18").Value = "LED1"
$ws.Range("B18").Value = "19-217/BHC-ZL1M2RY/3T"
$ws.Range("C18").Value = "LED0603"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "C72041"
$ws.Range("F18").Value = 0.0986

# --- View state ---------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I20").Select()
